# Appends three new paragraphs (Réplica 3 + its two body paragraphs) to the
# end of the document, right after the paragraph ending "... a cancelar."
# and before the section break.

$d = $word.ActiveDocument

# Collapsed range at the very end of the document body (end of the last
# paragraph, which currently ends with "a cancelar."). InsertXML on a
# collapsed range inserts content there without clobbering anything else.
$endRange = $d.Content
$endRange.Collapse(0)

$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-CR"/></w:rPr>'
$pPr = '<w:pPr><w:jc w:val="both"/>' + $rPr + '</w:pPr>'

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
$pPr
<w:r>$rPr<w:t>Réplica 3:</w:t></w:r>
</w:p>
<w:p>
$pPr
<w:r>$rPr<w:t>Hola, muchas gracias por ahondar más en el tema de las relaciones en matemáticas.</w:t></w:r>
</w:p>
<w:p>
$pPr
<w:r>$rPr<w:t xml:space="preserve">Otro de los ejemplos en donde podemos ver aplicado las relaciones en informática y que constantemente utilizamos (para redes sociales, para plataformas de pago, para servicios de </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r>$rPr<w:t>streaming</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>$rPr<w:t xml:space="preserve">, etc.) es en los sistemas en donde para navegar por una página web se debe estar registrado en ella. Aquí la relación es muy sencilla, el usuario o correo debe estar relacionado con una contraseña o </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:proofErr w:type="gramStart"/>
<w:r>$rPr<w:t>password</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:proofErr w:type="gramEnd"/>
<w:r>$rPr<w:t>, en caso de que no coincida con la registrada previamente se envía un error y no se le permitirá el ingreso al usuario.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$endRange.InsertXML($xml)

Write-Output "done"
